$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Artikelnummer"
$ws.Range("B1").Value = "Menge"
$ws.Range("C1").Value = "Empfänger"
$ws.Range("D1").Value = "Status"

# Empfänger column for rows 2-5 (establishes shared-string order: PL, HUN)
$ws.Range("C2").Value = "PL"
$ws.Range("C3").Value = "HUN"
$ws.Range("C4").Value = "HUN"
$ws.Range("C5").Value = "HUN"

# Artikelnummer column (KL2004..KL2009)
$ws.Range("A2").Value = "KL2004"
$ws.Range("A3").Value = "KL2005"
$ws.Range("A4").Value = "KL2006"
$ws.Range("A5").Value = "KL2007"
$ws.Range("A6").Value = "KL2008"
$ws.Range("A7").Value = "KL2009"

# Remaining Empfänger cells
$ws.Range("C6").Value = "3 HUN; 8 PL"
$ws.Range("C7").Value = "PL"

# Menge column
$ws.Range("B2").Value = 7
$ws.Range("B3").Value = 8
$ws.Range("B4").Value = 9
$ws.Range("B5").Value = 10
$ws.Range("B6").Value = 11
$ws.Range("B7").Value = 12

# Status column
$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("D5").Value = 1
$ws.Range("D6").Value = 0
$ws.Range("D7").Value = 0

# Final selection matches the target view state: active cell C6
$ws.Range("C6").Select()
